$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 173734
$ws.Range("C5").Value = 9855
$ws.Range("C6").Value = 504
$ws.Range("C7").Value = 5.67
